$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Main document body: remove the three demonstration paragraphs
#    ("A simple demonstration of a query :", the erroring "m:self."
#    field paragraph, and "End of demonstration.") leaving only the
#    trailing bookmark paragraph + sectPr untouched.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p3 = $d.Paragraphs.Item(3)
$bodyRng = $d.Range($p1.Range.Start, $p3.Range.End)
$bodyRng.Delete()

# ---------------------------------------------------------------------
# 2) Footer: collapse the first paragraph's runs (previously split up
#    with spell-check proofErr markers) into a single run reading
#    "A simple demonstration of a query :" tagged as English (US).
# ---------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$footerRng = $footer.Range.Duplicate
$found = $footerRng.Find.Execute("A simple demonstration of a query" + [char]160 + ":", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $footerRng.Text = "A simple demonstration of a query :"
    $footerRng.LanguageID = "en-US"
}

# ---------------------------------------------------------------------
# 3) Footer: reword the AQL error message.
# ---------------------------------------------------------------------
$footer2 = $sec.Footers.Item(1)
$errRng = $footer2.Range.Duplicate
$errRng.Find.Execute("Syntax error in AQL expression:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Invalid query statement:", 2) | Out-Null
